# Append two identical "Rest Day" workout log rows for Devi on 2025-02-09.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    $ws.Cells.Item($r, 1).Value = "Devi"

    # Column B holds a plain text date ("2025-02-09"), not a real date value.
    # A bare .Value assignment of a date-looking string gets auto-parsed into
    # a date serial by Excel, so force the cell to Text first, then drop the
    # now-unneeded explicit format back to Normal so the cell keeps the
    # workbook's default (unstyled) look, matching a plain data row.
    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2025-02-09"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 3).Value = "Rest Day"
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 0
}
